$d = $word.ActiveDocument

# Locate the target paragraph: "Using the higher order function filter(), define
# a function filter_long_words() ... that are longer than n." (exercise about
# filtering words longer than n).
$targetPara = $null
$nextPara = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $t = $para.Range.Text
    if (($t -like "*filter_long_words*") -and ($t -like "*Using the higher order function*")) {
        $targetPara = $para
        $nextPara = $d.Paragraphs.Item($i + 1)
        break
    }
}

# Highlight the whole exercise paragraph (all of its runs) in yellow.
$targetRange = $targetPara.Range
$targetRange.HighlightColorIndex = 7

# Move the "_GoBack" bookmark from the start of this paragraph to the start of
# the following paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $bm = $d.Bookmarks.Item("_GoBack")
    $bm.Delete()
}
$nextStart = $nextPara.Range.Start
$bmRange = $d.Range($nextStart, $nextStart)
$d.Bookmarks.Add("_GoBack", $bmRange)
